# Weekly data refresh: insert two new rows of Ciboulette price data at the
# top of the date-ordered block (rows 440-441), pushing the existing rows
# down by two (440->442, ..., 464->466).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 440; everything from 440 on shifts down.
$ws.Rows.Item(440).Insert()
$ws.Rows.Item(440).Insert()

# New row 440 - "Primera" quality
$ws.Range("A440").Value = 6
$ws.Range("B440").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C440").Value = "Metropolitana"
$ws.Range("D440").Value = "2022-07-11"
$ws.Range("D440").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E440").Value = 13
$ws.Range("F440").Value = 100112039
$ws.Range("G440").Value = "Ciboulette"
$ws.Range("H440").Value = "Sin especificar"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 230
$ws.Range("K440").Value = 2000
$ws.Range("L440").Value = 2000
$ws.Range("M440").Value = 2000
$ws.Range("N440").Value = "$/docena de atados"
$ws.Range("O440").Value = "Región Metropolitana"
$ws.Range("P440").Value = 667
$ws.Range("Q440").Value = 3
$ws.Range("R440").Value = "Hortaliza"

# New row 441 - "Segunda" quality
$ws.Range("A441").Value = 6
$ws.Range("B441").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C441").Value = "Metropolitana"
$ws.Range("D441").Value = "2022-07-11"
$ws.Range("D441").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E441").Value = 13
$ws.Range("F441").Value = 100112039
$ws.Range("G441").Value = "Ciboulette"
$ws.Range("H441").Value = "Sin especificar"
$ws.Range("I441").Value = "Segunda"
$ws.Range("J441").Value = 150
$ws.Range("K441").Value = 1500
$ws.Range("L441").Value = 1500
$ws.Range("M441").Value = 1500
$ws.Range("N441").Value = "$/docena de atados"
$ws.Range("O441").Value = "Región Metropolitana"
$ws.Range("P441").Value = 500
$ws.Range("Q441").Value = 3
$ws.Range("R441").Value = "Hortaliza"
